# Auto-update hourly job matches and history
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Posted-At column holds literal date-like text ("YYYY-MM-DD"), not real
# dates -- force text formatting first so Excel doesn't auto-coerce the
# assigned string into a date serial, then drop the temporary format so
# the cell keeps the sheet's (unstyled) default look.
$ws.Range("F2:F10").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Sr Field Engineer"
$ws.Range("B2").Value = "Striim"
$ws.Range("C2").Value = "Remote, US USA"
$ws.Range("D2").Value = 15.6
$ws.Range("E2").Value = "RAG, BigQuery, Kubernetes, Git, Snowflake, Databricks, BigQuery, Kafka, MySQL, Python"
$ws.Range("F2").Value = "2026-02-23"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=9fd06b655c3730ee"

# Row 3
$ws.Range("A3").Value = "Data Engineer II"
$ws.Range("B3").Value = "Corteva Agriscience"
$ws.Range("C3").Value = "Indianapolis, IN, US USA"
$ws.Range("D3").Value = 14.4
$ws.Range("E3").Value = "RAG, S3, EC2, FastAPI, Docker, Kubernetes, Git, Databricks, PySpark, Kafka"
$ws.Range("F3").Value = "2026-02-23"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=84bbb3c9d11dfa7f"

# Row 4
$ws.Range("A4").Value = "Senior Software Engineer"
$ws.Range("B4").Value = "project44"
$ws.Range("C4").Value = "Chicago, IL, US USA"
$ws.Range("D4").Value = 14.4
$ws.Range("E4").Value = "RAG, Copilot, Kinesis, Docker, Kubernetes, Git, Kafka, MongoDB, NoSQL, SQL"
$ws.Range("F4").Value = "2026-02-23"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=2ad7acbb70aad17a"

# Row 5
$ws.Range("A5").Value = "Software Engineer"
$ws.Range("B5").Value = "project44"
$ws.Range("C5").Value = "Chicago, IL, US USA"
$ws.Range("D5").Value = 14.4
$ws.Range("E5").Value = "RAG, Copilot, Kinesis, Docker, Kubernetes, Git, Kafka, MongoDB, NoSQL, SQL"
$ws.Range("F5").Value = "2026-02-23"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=166c9d26347adcb3"

# Row 6
$ws.Range("A6").Value = "Senior AI Engineer"
$ws.Range("B6").Value = "Medica Services Company LLC"
$ws.Range("C6").Value = "Minnetonka, MN, US USA"
$ws.Range("D6").Value = 12.2
$ws.Range("E6").Value = "AI Engineer, Generative AI, RAG, Copilot, Prompt Engineering, Docker, GitHub Actions, Terraform, Git, Python"
$ws.Range("F6").Value = "2026-02-23"
$ws.Range("G6").Value = "https://www.indeed.com/viewjob?jk=a79da7779f52d4b3"

# Row 7
$ws.Range("A7").Value = "Senior Software Engineer"
$ws.Range("B7").Value = "Kentan Staffing Solutions"
$ws.Range("C7").Value = "Melbourne, FL, US USA"
$ws.Range("D7").Value = 11.1
$ws.Range("E7").Value = "Kubernetes, CI/CD, Terraform, Git, MongoDB, NoSQL, SQL, R, Java, Scala"
$ws.Range("F7").Value = "2026-02-23"
$ws.Range("G7").Value = "https://www.indeed.com/viewjob?jk=ba7c1453cd1512f5"

# Row 8
$ws.Range("A8").Value = "Junior Software Engineer"
$ws.Range("B8").Value = "Kentan Staffing Solutions"
$ws.Range("C8").Value = "Melbourne, FL, US USA"
$ws.Range("D8").Value = 11.1
$ws.Range("E8").Value = "Docker, Kubernetes, CI/CD, Git, MongoDB, NoSQL, Python, SQL, R, Java"
$ws.Range("F8").Value = "2026-02-23"
$ws.Range("G8").Value = "https://www.indeed.com/viewjob?jk=726fbe9b2b7bd3f4"

# Row 9
$ws.Range("A9").Value = "Agentic AI Engineer-1"
$ws.Range("B9").Value = "Realign"
$ws.Range("C9").Value = "Boston, MA, US USA"
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = "AI Engineer, Generative AI, LangChain, RAG, Gemini, Copilot, Prompt Engineering, Python, R"
$ws.Range("F9").Value = "2026-02-23"
$ws.Range("G9").Value = "https://www.indeed.com/viewjob?jk=28b1ed79dd102956"

# Row 10
$ws.Range("A10").Value = "Perception Engineer - Data"
$ws.Range("B10").Value = "Forterra"
$ws.Range("C10").Value = "Arlington, VA, US USA"
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = "RAG, TensorFlow, PyTorch, Docker, Kubernetes, CI/CD, Python, R, Optimization"
$ws.Range("F10").Value = "2026-02-23"
$ws.Range("G10").Value = "https://www.indeed.com/viewjob?jk=24a4bdf67e7eef3b"

# Drop the scratch "@" text format now that the literal strings are safely
# stored -- restores the cells to the workbook's default (unstyled) look.
$ws.Range("F2:F10").ClearFormats()
